$d = $word.ActiveDocument
$hdr = $d.Sections(1).Headers.Item(1)
Write-Host "Header exists:" $hdr.Exists
$txt = $hdr.Range.Text
Write-Host "Header text:" ([System.String]::Join(",", [int[]][char[]]$txt))
